$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.413372394153498
$ws.Range("C2").Value = 0.201456015503112
$ws.Range("D2").Value = 0.1185488847114016
$ws.Range("E2").Value = 0.1291941946858399
$ws.Range("F2").Value = 1.803637194831531
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1737659992641758
$ws.Range("L2").Value = 0.272579169364171
$ws.Range("M2").Value = 0.3230976196828124
$ws.Range("O2").Value = 4.726910747837309
$ws.Range("B3").Value = 1.314378415876035
$ws.Range("C3").Value = 0.1887386274958658
$ws.Range("D3").Value = 0.1182101425572206
$ws.Range("E3").Value = 0.1301246217882692
$ws.Range("F3").Value = 1.816263475189238
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1750953665009884
$ws.Range("L3").Value = 0.2687900340461695
$ws.Range("M3").Value = 0.3071849845766366
$ws.Range("O3").Value = 4.765067934788078
$ws.Range("B4").Value = 1.253802580274055
$ws.Range("C4").Value = 0.1808928524815485
$ws.Range("D4").Value = 0.1180283284263695
$ws.Range("E4").Value = 0.1307297229449214
$ws.Range("F4").Value = 1.825028694161645
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1759581931711569
$ws.Range("L4").Value = 0.2665419847514272
$ws.Range("M4").Value = 0.2974856371997561
$ws.Range("O4").Value = 4.791342454178078
$ws.Range("B5").Value = 1.229171044146653
$ws.Range("C5").Value = 0.177686475829006
$ws.Range("D5").Value = 0.1179608494628575
$ws.Range("E5").Value = 0.1309848238460094
$ws.Range("F5").Value = 1.828855203617941
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1763215372722504
$ws.Range("L5").Value = 0.2656457351205148
$ws.Range("M5").Value = 0.2935512272236167
$ws.Range("O5").Value = 4.802764799976671
$ws.Range("B6").Value = 1.225084281260138
$ws.Range("C6").Value = 0.1771535123344705
$ws.Range("D6").Value = 0.117950044848893
$ws.Range("E6").Value = 0.1310276979949387
$ws.Range("F6").Value = 1.82950596984476
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1763825796478349
$ws.Range("L6").Value = 0.265498115628823
$ws.Range("M6").Value = 0.2928990261310318
$ws.Range("O6").Value = 4.804704661930742
$ws.Range("B7").Value = 1.253470171585718
$ws.Range("C7").Value = 0.1808496469739964
$ws.Range("D7").Value = 0.1180273915695231
$ws.Range("E7").Value = 0.1307331288198812
$ws.Range("F7").Value = 1.825079268927482
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1759630458151182
$ws.Range("L7").Value = 0.2665298171012935
$ws.Range("M7").Value = 0.2974325024971662
$ws.Range("O7").Value = 4.791493604314923
$ws.Range("B8").Value = 1.379197253359678
$ws.Range("C8").Value = 0.1970789102430217
$ws.Range("D8").Value = 0.1184266714329851
$ws.Range("E8").Value = 0.1295079962001919
$ws.Range("F8").Value = 1.807780549492783
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1742147073574989
$ws.Range("L8").Value = 0.2712564483201874
$ws.Range("M8").Value = 0.3175963537850208
$ws.Range("O8").Value = 4.73947647833603
$ws.Range("B9").Value = 1.627332369568535
$ws.Range("C9").Value = 0.2286013970677629
$ws.Range("D9").Value = 0.1194161389837944
$ws.Range("E9").Value = 0.1273731662755422
$ws.Range("F9").Value = 1.781894509591751
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1711549734841888
$ws.Range("L9").Value = 0.2811439992479166
$ws.Range("M9").Value = 0.3576914336833781
$ws.Range("O9").Value = 4.660069162567112
$ws.Range("B10").Value = 1.810546480302207
$ws.Range("C10").Value = 0.2515684342345708
$ws.Range("D10").Value = 0.1202675146623022
$ws.Range("E10").Value = 0.125966941002706
$ws.Range("F10").Value = 1.767779362188676
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1691304891481877
$ws.Range("L10").Value = 0.2887808919493011
$ws.Range("M10").Value = 0.387476564402391
$ws.Range("O10").Value = 4.615532086818973
$ws.Range("B11").Value = 1.894082098309525
$ws.Range("C11").Value = 0.261973451611226
$ws.Range("D11").Value = 0.1206815454333494
$ws.Range("E11").Value = 0.125362241172053
$ws.Range("F11").Value = 1.762423678517848
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1682577555639497
$ws.Range("L11").Value = 0.2923350893886152
$ws.Range("M11").Value = 0.4010956942395296
$ws.Range("O11").Value = 4.598274365031358
$ws.Range("B12").Value = 1.925740974398821
$ws.Range("C12").Value = 0.2659072357272407
$ws.Range("D12").Value = 0.1208421452284796
$ws.Range("E12").Value = 0.1251382753191099
$ws.Range("F12").Value = 1.760548882973524
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.167934187544823
$ws.Range("L12").Value = 0.293692399584998
$ws.Range("M12").Value = 0.4062626991946985
$ws.Range("O12").Value = 4.592171528068803
$ws.Range("B13").Value = 1.918921549962874
$ws.Range("C13").Value = 0.2650603116462946
$ws.Range("D13").Value = 0.1208073879653213
$ws.Range("E13").Value = 0.1251862872592024
$ws.Range("F13").Value = 1.760945834938951
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.168003566311727
$ws.Range("L13").Value = 0.2933995727195793
$ws.Range("M13").Value = 0.4051494638105169
$ws.Range("O13").Value = 4.593466648449748
$ws.Range("B14").Value = 1.896686186464422
$ws.Range("C14").Value = 0.2622972152873615
$ws.Range("D14").Value = 0.1206946817769534
$ws.Range("E14").Value = 0.1253437148095964
$ws.Range("F14").Value = 1.762266365391469
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1682309968867681
$ws.Range("L14").Value = 0.2924465280184165
$ws.Range("M14").Value = 0.4015205933759134
$ws.Range("O14").Value = 4.597763614164165
$ws.Range("B15").Value = 1.883069688813407
$ws.Range("C15").Value = 0.260603903063128
$ws.Range("D15").Value = 0.1206261419864845
$ws.Range("E15").Value = 0.1254407972385119
$ws.Range("F15").Value = 1.763095193309454
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.168371205042245
$ws.Range("L15").Value = 0.2918642436241612
$ws.Range("M15").Value = 0.399299064801788
$ws.Range("O15").Value = 4.600451942463565
$ws.Range("B16").Value = 1.805090924550484
$ws.Range("C16").Value = 0.2508875621264508
$ws.Range("D16").Value = 0.1202409922649039
$ws.Range("E16").Value = 0.1260071628043513
$ws.Range("F16").Value = 1.768150814089381
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1691884931710312
$ws.Range("L16").Value = 0.2885502198546135
$ws.Range("M16").Value = 0.386587900409225
$ws.Range("O16").Value = 4.616720382366964
$ws.Range("B17").Value = 1.757301111757101
$ws.Range("C17").Value = 0.2449157849360404
$ws.Range("D17").Value = 0.1200115421759236
$ws.Range("E17").Value = 0.1263635652952764
$ws.Range("F17").Value = 1.771525202052544
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1697022104330053
$ws.Range("L17").Value = 0.2865376194632745
$ws.Range("M17").Value = 0.3788076636251105
$ws.Range("O17").Value = 4.627469921995271
$ws.Range("B18").Value = 1.729831713786439
$ws.Range("C18").Value = 0.2414769605222489
$ws.Range("D18").Value = 0.1198820867812813
$ws.Range("E18").Value = 0.1265718535738491
$ws.Range("F18").Value = 1.773566333034552
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1700022261273304
$ws.Range("L18").Value = 0.2853875743668937
$ws.Range("M18").Value = 0.3743392572044968
$ws.Range("O18").Value = 4.633935342195883
$ws.Range("B19").Value = 1.720534196469487
$ws.Range("C19").Value = 0.2403119504864435
$ws.Range("D19").Value = 0.11983868881115
$ws.Range("E19").Value = 0.1266429426903799
$ws.Range("F19").Value = 1.774274644115792
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1701045863732968
$ws.Range("L19").Value = 0.2849994888654663
$ws.Range("M19").Value = 0.3728274705714298
$ws.Range("O19").Value = 4.636172937351887
$ws.Range("B20").Value = 1.76238656406349
$ws.Range("C20").Value = 0.2455519080328088
$ws.Range("D20").Value = 0.1200357071151998
$ws.Range("E20").Value = 0.1263252847094725
$ws.Range("F20").Value = 1.77115561425326
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1696470546887081
$ws.Range("L20").Value = 0.2867510836798033
$ws.Range("M20").Value = 0.379635204217692
$ws.Range("O20").Value = 4.626296366330337
$ws.Range("B21").Value = 1.903216565995479
$ws.Range("C21").Value = 0.2631089782539107
$ws.Range("D21").Value = 0.1207276830020234
$ws.Range("E21").Value = 0.125297338352627
$ws.Range("F21").Value = 1.76187433276877
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1681640074070962
$ws.Range("L21").Value = 0.2927261514897026
$ws.Range("M21").Value = 0.4025862186552516
$ws.Range("O21").Value = 4.596489754684825
$ws.Range("B22").Value = 1.995406447849859
$ws.Range("C22").Value = 0.2745463104228065
$ws.Range("D22").Value = 0.1212021544665163
$ws.Range("E22").Value = 0.1246547758332692
$ws.Range("F22").Value = 1.756701984709707
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1672350620801826
$ws.Range("L22").Value = 0.2966976751339843
$ws.Range("M22").Value = 0.4176426495161181
$ws.Range("O22").Value = 4.579529365649393
$ws.Range("B23").Value = 1.946189852775206
$ws.Range("C23").Value = 0.2684454678093857
$ws.Range("D23").Value = 0.1209468961151146
$ws.Range("E23").Value = 0.1249950502021986
$ws.Range("F23").Value = 1.759380779759553
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1677271743508868
$ws.Range("L23").Value = 0.294571952493726
$ws.Range("M23").Value = 0.4096016639101876
$ws.Range("O23").Value = 4.588350689299318
$ws.Range("B24").Value = 1.760087413501765
$ws.Range("C24").Value = 0.2452643341163707
$ws.Range("D24").Value = 0.1200247744873479
$ws.Range("E24").Value = 0.1263425808089753
$ws.Range("F24").Value = 1.771322389822913
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1696719760442171
$ws.Range("L24").Value = 0.2866545546155095
$ws.Range("M24").Value = 0.3792610589048664
$ws.Range("O24").Value = 4.62682604197434
$ws.Range("B25").Value = 1.560041525251904
$ws.Range("C25").Value = 0.220106987405444
$ws.Range("D25").Value = 0.1191265085268611
$ws.Range("E25").Value = 0.1279221371701809
$ws.Range("F25").Value = 1.788036576542254
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1719433703701938
$ws.Range("L25").Value = 0.2784034066239087
$ws.Range("M25").Value = 0.3467864406364356
$ws.Range("O25").Value = 4.679128858750005
